$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing date cell (A53) onto the
# new date cell (A54) so the new row's date uses the same style index as the
# rest of the date column, instead of Excel auto-creating a brand new style.
$ws.Cells.Item(53, 1).Copy()
$ws.Cells.Item(54, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 54 values.
$ws.Cells.Item(54, 1).Value = 45986
$ws.Cells.Item(54, 2).Value = 2025
$ws.Cells.Item(54, 3).Value = 1.049317648994741
$ws.Cells.Item(54, 4).Value = 2026
$ws.Cells.Item(54, 5).Value = 0.72625340902297
